$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended below the existing table (rows 1-13).
$ws.Range("A14").Value = 15
$ws.Range("B14").Value = 58
$ws.Range("A15").Value = 67
$ws.Range("A16").Value = 8

# Row 17 itself stays empty, but it must still exist so the sheet's
# used range (and therefore its <dimension>) extends down to row 17.
# Touching a format property on A17 registers the row/cell in the
# used range; clearing the format immediately after removes any
# visible formatting again, leaving row 17 present but blank.
$ws.Range("A17").NumberFormat = "General"
$ws.Range("A17").ClearFormats()
